# Shorten power plant lifetime to allow economic retirements pre-2030
$wb = $excel.ActiveWorkbook

# The "Lifetime" value for "All Plants" / "All years" lives on the MLfPPR sheet, cell B2.
$ws = $wb.Worksheets.Item("MLfPPR")

# Update the minimum lifetime for power plant retirements from 20 to 14 years.
$ws.Range("B2").Value = 14

# Make the MLfPPR sheet the active / selected sheet, matching the saved view state.
$ws.Activate()
$ws.Range("C2").Select()
